$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete earliest years (2005年-2009年), shifting the remaining
# rows (2010年 onward) up by five rows.
$ws.Range("A2:A6").EntireRow.Delete()

# After the shift, the sheet now ends at row 11 (2020年). Append a new row
# for 2021年, copying the row format/layout from the previous row (2020年)
# so blank columns keep the same look, then fill in this year's figures.
$ws.Range("A11:Q11").Copy($ws.Range("A12:Q12"))
$ws.Range("A12").Value = "2021年"
$ws.Range("N12").Value = 31773.1088999997
$ws.Range("O12").Value = 8199
